$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text values look numeric need an explicit text format so
# Excel does not coerce the assigned string into a Double (losing things
# like trailing zeros or the original formatting).
$textCells = @(
    "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D16", "D17", "D18", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D37", "D38", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D50", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.299.70"
$ws.Range("E2").Value = "  +0.67%  "

$ws.Range("D3").Value = "1.904.30"
$ws.Range("E3").Value = "  +0.54%  "

$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").Value = "306.56"
$ws.Range("E5").Value = "  -0.12%  "

$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("D7").Value = "0.5416"
$ws.Range("E7").Value = "  +3.99%  "

$ws.Range("D8").Value = "0.3810"
$ws.Range("E8").Value = "  +1.26%  "

$ws.Range("D9").Value = "0.07304"
$ws.Range("E9").Value = "  +0.39%  "

$ws.Range("D10").Value = "22.11"
$ws.Range("E10").Value = "  +4.48%  "

$ws.Range("D11").Value = "0.9025"
$ws.Range("E11").Value = "  +0.24%  "

$ws.Range("D12").Value = "0.08195"
$ws.Range("E12").Value = "  -0.22%  "

$ws.Range("D13").Value = "95.58"
$ws.Range("E13").Value = "  -0.85%  "

$ws.Range("D14").Value = "5.371"
$ws.Range("E14").Value = "  +1.11%  "

$ws.Range("D15").Value = "1.437.32"
$ws.Range("E15").Value = "  -24.62%  "

$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.16%  "

$ws.Range("D17").Value = "14.88"
$ws.Range("E17").Value = "  +1.98%  "

$ws.Range("D18").Value = "0.000008660"
$ws.Range("E18").Value = "  +0.53%  "

$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.26%  "

$ws.Range("D20").Value = "27.319.09"
$ws.Range("E20").Value = "  +0.63%  "

$ws.Range("D21").Value = "5.052"
$ws.Range("E21").Value = "  -0.58%  "

$ws.Range("D22").Value = "10.85"
$ws.Range("E22").Value = "  +1.30%  "

$ws.Range("D23").Value = "6.513"
$ws.Range("E23").Value = "  +1.44%  "

$ws.Range("B24").Value = "LidoDAOToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D24").Value = "2.317"
$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "148.88"
$ws.Range("E25").Value = "  +0.29%  "

$ws.Range("D26").Value = "18.39"
$ws.Range("E26").Value = "  +1.03%  "

$ws.Range("D27").Value = "1.747"
$ws.Range("E27").Value = "  +0.25%  "

$ws.Range("D28").Value = "116.53"
$ws.Range("E28").Value = "  +1.09%  "

$ws.Range("D29").Value = "4.853"
$ws.Range("E29").Value = "  +0.91%  "

$ws.Range("D30").Value = "4.672"
$ws.Range("E30").Value = "  -3.83%  "

$ws.Range("D31").Value = "0.09192"
$ws.Range("E31").Value = "  -0.16%  "

$ws.Range("D32").Value = "0.8267"
$ws.Range("E32").Value = "  +4.09%  "

$ws.Range("D33").Value = "0.05069"
$ws.Range("E33").Value = "  +0.91%  "

$ws.Range("E34").Value = "  +0.86%  "

$ws.Range("D35").Value = "3.021"
$ws.Range("E35").Value = "  +1.64%  "

$ws.Range("D36").Value = "3.316"
$ws.Range("E36").Value = "  -3.36%  "

$ws.Range("D37").Value = "2.690"
$ws.Range("E37").Value = "  +2.96%  "

$ws.Range("D38").Value = "0.6000"
$ws.Range("E38").Value = "  +4.79%  "

$ws.Range("E39").Value = "  -0.09%  "

$ws.Range("E40").Value = "  -0.03%  "

$ws.Range("D41").Value = "9.309"
$ws.Range("E41").Value = "  +3.23%  "

$ws.Range("D42").Value = "6.665"
$ws.Range("E42").Value = "  +1.55%  "

$ws.Range("D43").Value = "116.02"
$ws.Range("E43").Value = "  -0.40%  "

$ws.Range("D44").Value = "0.5132"
$ws.Range("E44").Value = "  +5.19%  "

$ws.Range("D45").Value = "0.1534"
$ws.Range("E45").Value = "  +1.29%  "

$ws.Range("D46").Value = "10.19"
$ws.Range("E46").Value = "  +1.04%  "

$ws.Range("D47").Value = "0.9998"
$ws.Range("E47").Value = "  -0.11%  "

$ws.Range("E48").Value = "  +0.97%  "

$ws.Range("D49").Value = "38.14"
$ws.Range("E49").Value = "  -0.45%  "

$ws.Range("D50").Value = "0.06097"
$ws.Range("E50").Value = "  +2.89%  "

$ws.Range("D51").Value = "63.34"
$ws.Range("E51").Value = "  -0.64%  "
